$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# I1: report date header 10/03/2023 -> 11/03/2023, preserve style/type (not a date re-parse)
$ws.Range("I1").Formula = '="11/03/2023"'
$ws.Range("I1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4163) | Out-Null
$excel.CutCopyMode = 0

# Row 2
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = 329.147
$ws.Range("D2").Value = 721.177
$ws.Range("E2").Value = 234
$ws.Range("F2").Value = 14
$ws.Range("G2").Value = 19
$ws.Range("H2").Value = 166
$ws.Range("I2").Value = 1444.3
$ws.Range("J2").Value = -50.06736827528906

# Row 3
$ws.Range("C3").Value = 121
$ws.Range("D3").Value = 122
$ws.Range("E3").Value = 1
$ws.Range("G3").Value = 1
$ws.Range("I3").Value = 165
$ws.Range("J3").Value = -26.06060606060606

# Row 4
$ws.Range("C4").Value = 347
$ws.Range("D4").Value = 348
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 11
$ws.Range("I4").Value = 57
$ws.Range("J4").Value = 510.5263157894737

# Row 5
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 687
$ws.Range("D5").Value = 960
$ws.Range("E5").Value = 198
$ws.Range("G5").Value = 6
$ws.Range("H5").Value = 59
$ws.Range("I5").Value = 1064
$ws.Range("J5").Value = -9.774436090225569

# Row 6
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 500
$ws.Range("D6").Value = 560
$ws.Range("E6").Value = 60
$ws.Range("F6").Value = 3
$ws.Range("G6").Value = 1
$ws.Range("I6").Value = 794
$ws.Range("J6").Value = -29.47103274559194

# Row 7
$ws.Range("C7").Value = 244
$ws.Range("D7").Value = 255
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0
$ws.Range("I7").Value = 264
$ws.Range("J7").Value = -3.409090909090906

# Row 8
$ws.Range("C8").Value = 58
$ws.Range("D8").Value = 78
$ws.Range("E8").Value = 19
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0
$ws.Range("I8").Value = 99
$ws.Range("J8").Value = -21.21212121212121

# Row 9
$ws.Range("B9").Value = 2
$ws.Range("C9").Value = 424
$ws.Range("D9").Value = 436
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = 2
$ws.Range("G9").Value = 10
$ws.Range("I9").Value = 1184
$ws.Range("J9").Value = -63.17567567567568

# Row 10
$ws.Range("D10").Value = 76
$ws.Range("E10").Value = 7
$ws.Range("G10").Value = 3
$ws.Range("I10").Value = 104
$ws.Range("J10").Value = -26.92307692307693

# Row 11
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("I11").Value = 3
$ws.Range("J11").Value = -100

# Row 12
$ws.Range("C12").Value = 40
$ws.Range("D12").Value = 44
$ws.Range("E12").Value = 4
$ws.Range("I12").Value = 83
$ws.Range("J12").Value = -46.98795180722891

